$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.723.40'
$ws.Range("E2").Value = '  +0.11%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.618.03'
$ws.Range("E3").Value = '  -0.66%  '
$ws.Range("E4").Value = '  +0.37%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '516.29'
$ws.Range("E5").Value = '  +1.77%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '154.51'
$ws.Range("E6").Value = '  -1.34%  '
$ws.Range("E7").Value = '  +0.11%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.590'
$ws.Range("E8").Value = '  +0.47%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.632.79'
$ws.Range("E9").Value = '  +0.10%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.70'
$ws.Range("E10").Value = '  +4.51%  '
$ws.Range("E11").Value = '  -0.03%  '
$ws.Range("E13").Value = '  +1.78%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.077.48'
$ws.Range("E14").Value = '  +0.89%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '60.747.14'
$ws.Range("E15").Value = '  +0.20%  '
$ws.Range("E16").Value = '  +0.20%  '
$ws.Range("E17").Value = '  +0.82%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.625.98'
$ws.Range("E18").Value = '  +0.11%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '357.21'
$ws.Range("E20").Value = '  +3.77%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.65'
$ws.Range("E21").Value = '  +2.22%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.21'
$ws.Range("E22").Value = '  +0.86%  '
$ws.Range("E23").Value = '  +0.19%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '60.88'
$ws.Range("E24").Value = '  +1.08%  '
$ws.Range("E25").Value = '  +1.00%  '
$ws.Range("E26").Value = '  +0.91%  '
$ws.Range("E27").Value = '  +0.88%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.994'
$ws.Range("E28").Value = '  +0.36%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0₃0846'
$ws.Range("E29").Value = '  -0.94%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.40'
$ws.Range("E30").Value = '  -1.95%  '
$ws.Range("E31").Value = '  +0.22%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '19.47'
$ws.Range("E32").Value = '  +0.52%  '
$ws.Range("E33").Value = '  +1.13%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '151.20'
$ws.Range("E34").Value = '  -3.30%  '
$ws.Range("E35").Value = '  +3.82%  '
$ws.Range("E36").Value = '  +0.73%  '
$ws.Range("E37").Value = '  -0.77%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.892'
$ws.Range("E38").Value = '  +6.83%  '
$ws.Range("E39").Value = '  +1.26%  '
$ws.Range("E40").Value = '  +0.21%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '36.44'
$ws.Range("E41").Value = '  +2.12%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.75'
$ws.Range("E42").Value = '  -0.70%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '292.77'
$ws.Range("E43").Value = '  -5.31%  '
$ws.Range("E44").Value = '  +0.06%  '
$ws.Range("E46").Value = '  -1.84%  '
$ws.Range("E47").Value = '  -0.34%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '19.81'
$ws.Range("E48").Value = '  -0.27%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '4.99'
$ws.Range("E49").Value = '  +2.34%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0237'
$ws.Range("E50").Value = '  +0.30%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '10.30'
$ws.Range("E51").Value = '  +0.30%  '
